# Apply the "Office Theme" colour scheme (the standard default Office
# palette) to the presentation's theme, replacing the current custom
# "Integral" / "Red Violet" colour scheme.
#
# PowerPoint stores the twelve theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) on the shared theme part used by the
# slide master; any Slide's ThemeColorScheme collection is a view onto
# that same shared theme, so editing it through slide 1 updates the
# whole deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# index : theme slot   : target RGB (Office default theme)
Set-ThemeColor $cs 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $cs 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $cs 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $cs 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $cs 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $cs 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $cs 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $cs 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $cs 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $cs 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $cs 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $cs 12 0x95 0x4F 0x72   # folHlink
